$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195, shifting existing rows 195:312 down to 196:313.
$ws.Rows.Item(195).Insert()

# Populate the new row 195 with the new record's data.
$ws.Cells.Item(195, 1).Value = 10
$ws.Cells.Item(195, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(195, 3).Value = "La Araucanía"
$ws.Cells.Item(195, 4).Value = 44830
$ws.Cells.Item(195, 5).Value = 9
$ws.Cells.Item(195, 6).Value = "Fruta"
$ws.Cells.Item(195, 7).Value = 100102
$ws.Cells.Item(195, 8).Value = "Cítricos"
$ws.Cells.Item(195, 9).Value = 100102006
$ws.Cells.Item(195, 10).Value = "Pomelo"
$ws.Cells.Item(195, 11).Value = "Start Ruby"
$ws.Cells.Item(195, 12).Value = "Primera"
$ws.Cells.Item(195, 13).Value = 180
$ws.Cells.Item(195, 14).Value = 15000
$ws.Cells.Item(195, 15).Value = 16000
$ws.Cells.Item(195, 16).Value = 15556
$ws.Cells.Item(195, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(195, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(195, 19).Value = 1037
$ws.Cells.Item(195, 20).Value = 15
